# Injects the template3.docx "dev injector mvp" content: twenty-one
# paragraphs (plain text paragraphs, two MERGEFIELD //@blockNNNN fields
# built from raw fldChar/instrText runs, a proofErr-wrapped paragraph,
# and the original _GoBack bookmark now anchored to "Десятый абзац").
#
# We replace the whole body story (both original paragraphs) in one shot
# via Range.InsertXML so the emitted run/field structure is byte-for-byte
# what Word would have produced, rather than relying on Find/Replace or
# Fields.Add (which synthesizes a different, collapsed field encoding).

$d = $word.ActiveDocument

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
<w:p><w:r><w:t>Нулевой абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Первый абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Второй абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Третий абзац</w:t></w:r></w:p>
<w:p>
<w:r><w:fldChar w:fldCharType="begin"/></w:r>
<w:r><w:instrText xml:space="preserve"> MERGEFIELD  //@block1788</w:instrText></w:r>
<w:r><w:instrText xml:space="preserve"> \* MERGEFORMAT </w:instrText></w:r>
<w:r><w:fldChar w:fldCharType="separate"/></w:r>
<w:r><w:rPr><w:noProof/></w:rPr><w:t>«//@block</w:t></w:r>
<w:r><w:rPr><w:noProof/></w:rPr><w:t>1788</w:t></w:r>
<w:r><w:rPr><w:noProof/></w:rPr><w:t>»</w:t></w:r>
<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>
</w:p>
<w:p><w:r><w:t>Пятый абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Шестой абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Седьмой абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Восьмой абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Девятый абзац</w:t></w:r></w:p>
<w:p>
<w:r><w:t>Десятый абзац</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Одинадцатый</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> абзац</w:t></w:r>
</w:p>
<w:p><w:r><w:t>Двенадцатый абзац</w:t></w:r></w:p>
<w:p>
<w:r><w:fldChar w:fldCharType="begin"/></w:r>
<w:r><w:instrText xml:space="preserve"> MERGEFIELD  //@block1789 \* MERGEFORMAT </w:instrText></w:r>
<w:r><w:fldChar w:fldCharType="separate"/></w:r>
<w:r><w:rPr><w:noProof/></w:rPr><w:t>«//@block</w:t></w:r>
<w:r><w:rPr><w:noProof/></w:rPr><w:t>1789</w:t></w:r>
<w:r><w:rPr><w:noProof/></w:rPr><w:t>»</w:t></w:r>
<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>
</w:p>
<w:p><w:r><w:t>Четырнадцатый абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Пятнадцатый абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Шестнадцатый абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Семнадцатый абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Восемнадцатый абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Девятнадцатый абзац</w:t></w:r></w:p>
<w:p><w:r><w:t>Двадцатый абзац</w:t></w:r></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

# Range.InsertXML replaces the target range's contents in place; the
# document's sectPr lives outside $d.Content so it (and its page setup)
# is left untouched.
$d.Content.InsertXML($xml)
